# Append 8 new coded-segment rows (147-154) to Sheet1, mirroring the
# formatting of the last existing row (146) and filling in the new
# "Event month" / "Event year" segment data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = $ws.Range("A146:M146")

# Copy row 146 (values + formatting) down into each new row first so that
# styles/number formats match the rest of the table, then overwrite the
# cells that actually differ per row.
for ($r = 147; $r -le 154; $r++) {
    $dst = $ws.Range("A" + $r + ":M" + $r)
    $srcRow.Copy($dst)
}

# Helper: assign a value that Excel would otherwise auto-coerce to a number
# (e.g. "2008") while preserving the destination cell's existing style.
# We temporarily borrow a Text-formatted cell's number format, assign the
# value (which now stays a string), then paste the original column's
# format back on top - this avoids leaving a quotePrefix/new style behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $ws.Cells.Item(146, 6).Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $val
    $ws.Cells.Item(146, $col).Copy()
    $cell.PasteSpecial(-4122)
}

# Row 147
$ws.Cells.Item(147, 5).Value = "Event month"
$ws.Cells.Item(147, 6).Value = "2: 1538"
$ws.Cells.Item(147, 7).Value = "2: 1546"
$ws.Cells.Item(147, 9).Value = "September"
$ws.Cells.Item(147, 10).Value = 9
$ws.Cells.Item(147, 11).Value = 0.030062
$ws.Cells.Item(147, 13).Value = "11/14/18 13:28:00"

# Row 148
$ws.Cells.Item(148, 5).Value = "Event month"
$ws.Cells.Item(148, 6).Value = "2: 1557"
$ws.Cells.Item(148, 7).Value = "2: 1561"
$ws.Cells.Item(148, 9).Value = "March"
$ws.Cells.Item(148, 10).Value = 5
$ws.Cells.Item(148, 11).Value = 0.016701
$ws.Cells.Item(148, 13).Value = "11/14/18 13:29:00"

# Row 149
$ws.Cells.Item(149, 5).Value = "Event year"
$ws.Cells.Item(149, 6).Value = "2: 1548"
$ws.Cells.Item(149, 7).Value = "2: 1551"
Set-TextValue 149 9 "2008"
$ws.Cells.Item(149, 10).Value = 4
$ws.Cells.Item(149, 11).Value = 0.013361
$ws.Cells.Item(149, 13).Value = "11/14/18 13:29:00"

# Row 150
$ws.Cells.Item(150, 5).Value = "Event year"
$ws.Cells.Item(150, 6).Value = "2: 1563"
$ws.Cells.Item(150, 7).Value = "2: 1566"
Set-TextValue 150 9 "2010"
$ws.Cells.Item(150, 10).Value = 4
$ws.Cells.Item(150, 11).Value = 0.013361
$ws.Cells.Item(150, 13).Value = "11/14/18 13:29:00"

# Row 151
$ws.Cells.Item(151, 5).Value = "I"
$ws.Cells.Item(151, 6).Value = "2: 1538"
$ws.Cells.Item(151, 7).Value = "2: 1546"
$ws.Cells.Item(151, 9).Value = "September"
$ws.Cells.Item(151, 10).Value = 9
$ws.Cells.Item(151, 11).Value = 0.030062
$ws.Cells.Item(151, 13).Value = "11/14/18 13:29:00"

# Row 152
$ws.Cells.Item(152, 5).Value = "I"
$ws.Cells.Item(152, 6).Value = "2: 1548"
$ws.Cells.Item(152, 7).Value = "2: 1551"
Set-TextValue 152 9 "2008"
$ws.Cells.Item(152, 10).Value = 4
$ws.Cells.Item(152, 11).Value = 0.013361
$ws.Cells.Item(152, 13).Value = "11/14/18 13:29:00"

# Row 153
$ws.Cells.Item(153, 5).Value = "J"
$ws.Cells.Item(153, 6).Value = "2: 1557"
$ws.Cells.Item(153, 7).Value = "2: 1561"
$ws.Cells.Item(153, 9).Value = "March"
$ws.Cells.Item(153, 10).Value = 5
$ws.Cells.Item(153, 11).Value = 0.016701
$ws.Cells.Item(153, 13).Value = "11/14/18 13:29:00"

# Row 154
$ws.Cells.Item(154, 5).Value = "J"
$ws.Cells.Item(154, 6).Value = "2: 1563"
$ws.Cells.Item(154, 7).Value = "2: 1566"
Set-TextValue 154 9 "2010"
$ws.Cells.Item(154, 10).Value = 4
$ws.Cells.Item(154, 11).Value = 0.013361
$ws.Cells.Item(154, 13).Value = "11/14/18 13:29:00"
